# [Billing] edit excel template
# Rebuilds the invoice template sheet: adds a "Company ID / Period / INVOICES"
# summary header above the existing table, moves the table down to rows 7-8,
# adds a "Summary" total row, and updates the jxls comments accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Capture the existing (pre-edit) header/detail rows before we overwrite
#    anything, since the new layout re-uses the same text in new rows.
# ---------------------------------------------------------------------------
$headerVals = @()
for ($c = 1; $c -le 7; $c++) {
    $headerVals += $ws.Cells.Item(3, $c).Text
}
$detailVals = @()
for ($c = 1; $c -le 6; $c++) {
    $detailVals += $ws.Cells.Item(4, $c).Text
}

# ---------------------------------------------------------------------------
# 2. Clear the old A3:G4 content - it gets rebuilt at A7:G8.
# ---------------------------------------------------------------------------
$ws.Range("A3:G4").Clear()

# ---------------------------------------------------------------------------
# 3. New summary rows 3-5 (Company ID / Period / INVOICES banner)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Company ID:"
$ws.Range("B3").Value = '${companyId}'

$ws.Range("A4").Value = "Period:"
$ws.Range("B4").Formula = '=TEXT(MIN(B8),"mmmm dd, yyyy") & " - " & TEXT(MAX(B8), "mmmm dd, yyyy")'
$ws.Range("B4").NumberFormat = "[$-409]mmmm\ d\,\ yyyy;@"

$ws.Range("A3:A4").Font.Bold = $true
$ws.Range("A3:A4").HorizontalAlignment = -4152   # xlRight

$ws.Range("A5:G5").Merge()
$ws.Range("A5").Value = "INVOICES"
$ws.Range("A5:G5").Font.Bold = $true
$ws.Range("A5:G5").Font.Size = 14
$ws.Range("A5:G5").HorizontalAlignment = -4108   # xlCenter
$ws.Rows.Item(5).RowHeight = 19

$ws.Rows.Item(6).RowHeight = 16

# ---------------------------------------------------------------------------
# 4. Header row (was row 3) rebuilt at row 7 with shading + borders
# ---------------------------------------------------------------------------
for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item(7, $c).Value = $headerVals[$c - 1]
}
$hdr = $ws.Range("A7:G7")
$hdr.Font.Bold = $true
$hdr.Interior.Pattern = -4124       # xlPatternSolid
$hdr.Interior.ThemeColor = 5        # matches theme="4" (Accent1)
$hdr.Interior.TintAndShade = 0.6
$ws.Range("E7:F7").HorizontalAlignment = -4152   # xlRight
$ws.Range("G7").HorizontalAlignment = -4152      # xlRight

$hdr.Borders.Item(8).LineStyle = 1       # xlEdgeTop
$hdr.Borders.Item(8).Weight = 4          # xlMedium
$hdr.Borders.Item(9).LineStyle = 1       # xlEdgeBottom
$hdr.Borders.Item(9).Weight = 2          # xlThin
$ws.Range("A7").Borders.Item(7).LineStyle = 1
$ws.Range("A7").Borders.Item(7).Weight = 4
$ws.Range("G7").Borders.Item(10).LineStyle = 1
$ws.Range("G7").Borders.Item(10).Weight = 4
$hdr.Borders.Item(11).LineStyle = 1      # xlInsideVertical
$hdr.Borders.Item(11).Weight = 2

# ---------------------------------------------------------------------------
# 5. Detail row (was row 4) rebuilt at row 8, G8 recomputed as a (now-bad)
#    formula E8*F8 - both operands are still template placeholders so this
#    legitimately evaluates to #VALUE!, matching the authored template.
# ---------------------------------------------------------------------------
for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(8, $c).Value = $detailVals[$c - 1]
}
$ws.Range("G8").Formula = "=E8*F8"
$ws.Range("B8:C8").NumberFormat = "[$-409]mmmm\ d\,\ yyyy;@"

$det = $ws.Range("A8:G8")
$det.Borders.Item(8).LineStyle = 1       # xlEdgeTop
$det.Borders.Item(8).Weight = 2          # xlThin
$ws.Range("A8").Borders.Item(7).LineStyle = 1
$ws.Range("A8").Borders.Item(7).Weight = 4
$ws.Range("G8").Borders.Item(10).LineStyle = 1
$ws.Range("G8").Borders.Item(10).Weight = 4
$det.Borders.Item(11).LineStyle = 1      # xlInsideVertical
$det.Borders.Item(11).Weight = 2

$ws.Rows.Item(8).RowHeight = 16

# ---------------------------------------------------------------------------
# 6. Summary row 9 - label + total
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "Summary:"
$ws.Range("A9:F9").Font.Bold = $true
$sumRow = $ws.Range("A9:F9")
$sumRow.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$sumRow.Borders.Item(8).Weight = 4       # xlMedium

$ws.Range("G9").Formula = "=SUM(G8)"
$ws.Range("G9").Font.Bold = $true
$ws.Range("G9").HorizontalAlignment = -4152   # xlRight
$ws.Range("G9").Borders.Item(8).LineStyle = 1
$ws.Range("G9").Borders.Item(8).Weight = 4

# ---------------------------------------------------------------------------
# 7. Column widths (characters) and selection / zoom
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 28.4167
$ws.Columns.Item(2).ColumnWidth = 22.25
$ws.Columns.Item(3).ColumnWidth = 26.25
$ws.Columns.Item(4).ColumnWidth = 20.25
$ws.Columns.Item(5).ColumnWidth = 24.9167
$ws.Columns.Item(6).ColumnWidth = 27.0833
$ws.Columns.Item(7).ColumnWidth = 16.4167

$ws.Range("A5:G5").Select()
$excel.ActiveWindow.Zoom = 150

# ---------------------------------------------------------------------------
# 8. Comments - update jxls directives and re-anchor the "each" comment to
#    the new detail row.
# ---------------------------------------------------------------------------
$ws.Range("A1").Comment.Text('jx:area(lastCell="G9")')

$ws.Range("A4").Comment.Delete()
$ws.Range("A8").AddComment("jx:each(items=`"invoices`" var=`"invoice`" lastCell=`"G8`")" + [char]10)
